# Auto-generated script to apply market-price data refresh to Halicarnassus_Profits sheets
# Updates computed price/profit columns (H:N) for specific rows across multiple worksheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 131.125
$ws.Range("I33").Value = 108.166664
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 108.166664
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 120.833336
$ws.Range("N33").Value = -658
$ws.Range("H70").Value = 1199.5714
$ws.Range("I70").Value = 945
$ws.Range("J70").Value = 1301.4
$ws.Range("K70").Value = 2835
$ws.Range("L70").Value = 3904.2
$ws.Range("M70").Value = -2565
$ws.Range("N70").Value = -4444.200000000001
$ws.Range("H73").Value = 1199.5714
$ws.Range("I73").Value = 945
$ws.Range("J73").Value = 1301.4
$ws.Range("K73").Value = 2835
$ws.Range("L73").Value = 3904.2
$ws.Range("M73").Value = -1899
$ws.Range("N73").Value = -5776.200000000001
$ws.Range("H115").Value = 819.75
$ws.Range("I115").Value = 1028.3334
$ws.Range("J115").Value = 194
$ws.Range("K115").Value = 3085.0002
$ws.Range("L115").Value = 582
$ws.Range("M115").Value = -1518.0002
$ws.Range("N115").Value = -3716
$ws.Range("H116").Value = 7359.1665
$ws.Range("I116").Value = 13187.5
$ws.Range("J116").Value = 4445
$ws.Range("K116").Value = 13187.5
$ws.Range("L116").Value = 4445
$ws.Range("M116").Value = -9745.5
$ws.Range("N116").Value = -11329
$ws.Range("H125").Value = 3508.0588
$ws.Range("I125").Value = 3456
$ws.Range("J125").Value = 3633
$ws.Range("K125").Value = 31104
$ws.Range("L125").Value = 32697
$ws.Range("M125").Value = -28644
$ws.Range("N125").Value = -37617
$ws.Range("H132").Value = 4609.2
$ws.Range("I132").Value = 1586.2667
$ws.Range("J132").Value = 13678
$ws.Range("K132").Value = 4758.800099999999
$ws.Range("L132").Value = 41034
$ws.Range("M132").Value = -2228.800099999999
$ws.Range("N132").Value = -46094
$ws.Range("H137").Value = 2214.923
$ws.Range("I137").Value = 1383.7693
$ws.Range("J137").Value = 3046.077
$ws.Range("K137").Value = 4151.3079
$ws.Range("L137").Value = 9138.231
$ws.Range("M137").Value = -1601.3079

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3729.0605
$ws.Range("I32").Value = 864.10345
$ws.Range("J32").Value = 24500
$ws.Range("K32").Value = 864.10345
$ws.Range("L32").Value = 24500
$ws.Range("M32").Value = -577.10345
$ws.Range("N32").Value = -25074
$ws.Range("H45").Value = 2112.923
$ws.Range("I45").Value = 1572.625
$ws.Range("J45").Value = 2977.4
$ws.Range("K45").Value = 1572.625
$ws.Range("L45").Value = 2977.4
$ws.Range("M45").Value = -1195.625
$ws.Range("N45").Value = -3731.4
$ws.Range("H132").Value = 2923.5334
$ws.Range("I132").Value = 2917.9092
$ws.Range("J132").Value = 2939
$ws.Range("K132").Value = 8753.7276
$ws.Range("L132").Value = 8817
$ws.Range("M132").Value = -6223.7276
$ws.Range("N132").Value = -13877
$ws.Range("H139").Value = 35000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 7387.5
$ws.Range("I25").Value = 550
$ws.Range("J25").Value = 9666.666999999999
$ws.Range("K25").Value = 550
$ws.Range("L25").Value = 9666.666999999999
$ws.Range("M25").Value = -315
$ws.Range("N25").Value = -10136.667
$ws.Range("H138").Value = 59000
$ws.Range("I138").Value = 25000
$ws.Range("J138").Value = 149666.67
$ws.Range("K138").Value = 25000
$ws.Range("L138").Value = 149666.67
$ws.Range("M138").Value = -19860
$ws.Range("N138").Value = -159946.67

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 61.35
$ws.Range("I7").Value = 45.23077
$ws.Range("J7").Value = 91.28570999999999
$ws.Range("K7").Value = 45.23077
$ws.Range("L7").Value = 91.28570999999999
$ws.Range("M7").Value = 67.76922999999999
$ws.Range("N7").Value = -317.28571
$ws.Range("H22").Value = 1050.2858
$ws.Range("I22").Value = 926.5
$ws.Range("J22").Value = 1143.125
$ws.Range("K22").Value = 926.5
$ws.Range("L22").Value = 1143.125
$ws.Range("M22").Value = -576.5
$ws.Range("N22").Value = -1843.125
$ws.Range("H31").Value = 5503.698
$ws.Range("I31").Value = 4251.7856
$ws.Range("J31").Value = 5953.1025
$ws.Range("K31").Value = 4251.7856
$ws.Range("L31").Value = 5953.1025
$ws.Range("M31").Value = -3956.7856
$ws.Range("N31").Value = -6543.1025
$ws.Range("H34").Value = 5503.698
$ws.Range("I34").Value = 4251.7856
$ws.Range("J34").Value = 5953.1025
$ws.Range("K34").Value = 4251.7856
$ws.Range("L34").Value = 5953.1025
$ws.Range("M34").Value = -4049.7856
$ws.Range("N34").Value = -6357.1025
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H58").Value = 5097.8
$ws.Range("I58").Value = 2499.5
$ws.Range("J58").Value = 6830
$ws.Range("K58").Value = 2499.5
$ws.Range("L58").Value = 6830
$ws.Range("M58").Value = -2296.5
$ws.Range("N58").Value = -7236
$ws.Range("H105").Value = 1348.25
$ws.Range("I105").Value = 798
$ws.Range("J105").Value = 2999
$ws.Range("K105").Value = 798
$ws.Range("L105").Value = 2999
$ws.Range("M105").Value = 949
$ws.Range("N105").Value = -6493
$ws.Range("H136").Value = 5097.8
$ws.Range("I136").Value = 2499.5
$ws.Range("J136").Value = 6830
$ws.Range("K136").Value = 7498.5
$ws.Range("L136").Value = 20490
$ws.Range("M136").Value = -4948.5
$ws.Range("N136").Value = -25590
$ws.Range("N48").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2500
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 2500
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 7500
$ws.Range("N41").Value = -8176
$ws.Range("H107").Value = 411.36365
$ws.Range("I107").Value = 302.4
$ws.Range("J107").Value = 443.41177
$ws.Range("K107").Value = 907.1999999999999
$ws.Range("L107").Value = 1330.23531
$ws.Range("M107").Value = 1012.8
$ws.Range("N107").Value = -5170.23531
$ws.Range("H137").Value = 3966.5
$ws.Range("I137").Value = 3999.5
$ws.Range("J137").Value = 3950
$ws.Range("K137").Value = 11998.5
$ws.Range("L137").Value = 11850
$ws.Range("M137").Value = -6898.5
$ws.Range("N137").Value = -22050
$ws.Range("H138").Value = 5302.375
$ws.Range("I138").Value = 2532.25
$ws.Range("J138").Value = 8072.5
$ws.Range("K138").Value = 7596.75
$ws.Range("L138").Value = 24217.5
$ws.Range("M138").Value = -2456.75
$ws.Range("N138").Value = -34497.5
$ws.Range("H139").Value = 857.5
$ws.Range("I139").Value = 857.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2572.5
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2567.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 118
$ws.Range("I9").Value = 123.6
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 123.6
$ws.Range("L9").Value = 123.6
$ws.Range("M9").Value = 46.40000000000001
$ws.Range("N9").Value = -430
$ws.Range("H97").Value = 605.6
$ws.Range("I97").Value = 678.36365
$ws.Range("J97").Value = 405.5
$ws.Range("K97").Value = 678.36365
$ws.Range("L97").Value = 405.5
$ws.Range("M97").Value = -182.36365
$ws.Range("N97").Value = -1397.5
$ws.Range("H132").Value = 96672.27
$ws.Range("I132").Value = 130524.375
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 391573.125
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -389043.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 20000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20590
$ws.Range("H122").Value = 4199
$ws.Range("I122").Value = 2949
$ws.Range("J122").Value = 4865.6665
$ws.Range("K122").Value = 8847
$ws.Range("L122").Value = 14596.9995
$ws.Range("M122").Value = -6397
$ws.Range("M29").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H107").Value = 814
$ws.Range("I107").Value = 814
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2442
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -522
$ws.Range("H132").Value = 2370.5715
$ws.Range("I132").Value = 2182.3333
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6546.999899999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -4016.999899999999
$ws.Range("N132").Value = -15560
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()
